$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values could be misread as numbers -> force text format first
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '62.701.45'
$ws.Range("E2").Value = '  +1.00%  '

# Row 3
$ws.Range("D3").Value = '2.437.74'
$ws.Range("E3").Value = '  +1.16%  '

# Row 4
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").Value = '566.26'
$ws.Range("E5").Value = '  +0.80%  '

# Row 6
$ws.Range("D6").Value = '145.25'
$ws.Range("E6").Value = '  +2.02%  '

# Row 7
$ws.Range("E7").Value = '  -0.08%  '

# Row 8
$ws.Range("E8").Value = '  +1.26%  '

# Row 9
$ws.Range("E9").Value = '  +2.03%  '

# Row 10
$ws.Range("E10").Value = '  +0.59%  '

# Row 11
$ws.Range("D11").Value = '5.25'
$ws.Range("E11").Value = '  -0.78%  '

# Row 12
$ws.Range("E12").Value = '  +1.40%  '

# Row 13
$ws.Range("D13").Value = '26.81'
$ws.Range("E13").Value = '  +5.03%  '

# Row 14
$ws.Range("E14").Value = '  +6.52%  '

# Row 15
$ws.Range("E15").Value = '  +1.14%  '

# Row 16
$ws.Range("D16").Value = '62.468.85'
$ws.Range("E16").Value = '  +0.59%  '

# Row 17
$ws.Range("D17").Value = '2.435.42'
$ws.Range("E17").Value = '  +1.06%  '

# Row 18
$ws.Range("E18").Value = '  -0.22%  '

# Row 19
$ws.Range("E19").Value = '  +1.47%  '

# Row 20
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").Value = '4.18'
$ws.Range("E20").Value = '  +0.94%  '

# Row 21
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = '323.54'
$ws.Range("E21").Value = '  +0.98%  '

# Row 22
$ws.Range("E22").Value = '  -0.06%  '

# Row 23
$ws.Range("D23").Value = '67.23'
$ws.Range("E23").Value = '  +1.87%  '

# Row 24
$ws.Range("E24").Value = '  +3.35%  '

# Row 25
$ws.Range("D25").Value = '8.73'
$ws.Range("E25").Value = '  -1.12%  '

# Row 26
$ws.Range("E26").Value = '  +8.75%  '

# Row 27
$ws.Range("D27").Value = '566.28'

# Row 28
$ws.Range("D28").Value = '2.557.06'
$ws.Range("E28").Value = '  +1.15%  '

# Row 29
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.14%  '

# Row 30
$ws.Range("D30").Value = '8.41'
$ws.Range("E30").Value = '  +3.26%  '

# Row 31
$ws.Range("D31").Value = '1.45'
$ws.Range("E31").Value = '  +3.15%  '

# Row 32
$ws.Range("E32").Value = '  -0.10%  '

# Row 33
$ws.Range("E33").Value = '  +0.44%  '

# Row 34
$ws.Range("D34").Value = '1.55'
$ws.Range("E34").Value = '  +1.67%  '

# Row 35
$ws.Range("D35").Value = '4.88'
$ws.Range("E35").Value = '  +4.52%  '

# Row 36
$ws.Range("E36").Value = '  -0.10%  '

# Row 37
$ws.Range("E37").Value = '  +1.02%  '

# Row 38
$ws.Range("D38").Value = '5.44'
$ws.Range("E38").Value = '  -0.46%  '

# Row 39
$ws.Range("D39").Value = '18.78'
$ws.Range("E39").Value = '  +1.01%  '

# Row 40
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").Value = '148.67'
$ws.Range("E40").Value = '  -1.96%  '

# Row 41
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = '1.83'
$ws.Range("E41").Value = '  +2.31%  '

# Row 42
$ws.Range("E42").Value = '  +0.41%  '

# Row 43
$ws.Range("D43").Value = '2.42'
$ws.Range("E43").Value = '  +6.64%  '

# Row 44
$ws.Range("D44").Value = '148.66'
$ws.Range("E44").Value = '  +0.66%  '

# Row 45
$ws.Range("E45").Value = '  +1.62%  '

# Row 46
$ws.Range("D46").Value = '0.0537'
$ws.Range("E46").Value = '  +1.23%  '

# Row 47
$ws.Range("E47").Value = '  +3.75%  '

# Row 48
$ws.Range("E48").Value = '  +1.59%  '

# Row 49
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = '0.0232'
$ws.Range("E49").Value = '  +3.19%  '

# Row 50
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").Value = '0.0928'
$ws.Range("E50").Value = '  +1.42%  '

# Row 51
$ws.Range("E51").Value = '  +0.56%  '
